$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values shift from 46060 (2026-02-07) to 46061 (2026-02-08)
# for every data row, from row 2 through row 204.
$ws.Range("C2:C204").Value = 46061
